# Final man. annotated files
# Apply the COP25 Malaysia annotation edits:
#  - remove the "Time" column (the "Scale" column is repurposed into "Shape")
#  - the "Relevance" column (B) switches from yes/no text to a 0/1/2 numeric code
#  - rows 15-17 get new Topic/Unit/Shape/Principle/explanation annotations

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column structure -------------------------------------------------
# Before: A=Unit(unit col) ... D=Unit, E=Scale, F=Time, G=Principle, H=30 word explanation, I=Notes
# Delete the "Time" column (column F) -- everything right of it shifts one left.
$ws.Columns.Item(6).Delete()

# The old "Scale" column (now column E) is renamed to "Shape".
$ws.Range("E1").Value = "Shape"

# --- Relevance column (B): yes/no text -> numeric code -----------------
$ws.Range("B2").Value = 0
$ws.Range("B3").Value = 0
$ws.Range("B4").Value = 0
$ws.Range("B5").Value = 0
$ws.Range("B6").Value = 0
$ws.Range("B7").Value = 1
$ws.Range("B8").Value = 0
$ws.Range("B9").Value = 1
$ws.Range("B10").Value = 0
$ws.Range("B11").Value = 0
$ws.Range("B12").Value = 0
$ws.Range("B13").Value = 0
$ws.Range("B14").Value = 0
$ws.Range("B15").Value = 2
$ws.Range("B16").Value = 2
$ws.Range("B17").Value = 2
$ws.Range("B18").Value = 1
$ws.Range("B19").Value = 1
$ws.Range("B20").Value = 0

# --- Row 15 annotation ---------------------------------------------------
$ws.Range("C15").Value = "UNFCCC agreements and principles"
$ws.Range("D15").Value = "responsibility"
$ws.Range("E15").Value = "n.a."
$ws.Range("F15").Value = "egalitarian"
$ws.Range("G15").Value = "Prescribing the moral obligation of developed countries to take the lead, motivated by the foundations of the convention, these are egalitarian. "

# --- Row 16 annotation ---------------------------------------------------
$ws.Range("C16").Value = "UNFCCC agreements and principles"
$ws.Range("D16").Value = "responsibility"
$ws.Range("E16").Value = "n.a."
$ws.Range("F16").Value = "egalitarian"
$ws.Range("G16").Value = "Value judgement on the need to take on action by developed countries. Egalitarian motviation by means of the focus on developed countries to contribute. "

# --- Row 17 annotation ---------------------------------------------------
$ws.Range("C17").Value = "GCF"
$ws.Range("D17").Value = "financial resources"
$ws.Range("E17").Value = "n.a."
$ws.Range("F17").Value = "prioritarian"
$ws.Range("G17").Value = "Value judgement on the need for redistribution of finances. Motivated by the prioritarian idea to help the worst off. "

# --- Selection / view state ----------------------------------------------
$ws.Range("B8").Select()
